# Update the cached "date last saved/printed" auto-field text on the
# Slide Master and on every Slide Layout's Date placeholder, from
# 10/23/24 to 10/24/24 (the field's id/type stay {48A87A34-...}
# / datetimeFigureOut; only the literal cached text moves forward a day).

$OLD_DATE = "10/23/24"
$NEW_DATE = "10/24/24"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)

        $isDatePlaceholder = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }

        if ($isDatePlaceholder -and $sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $OLD_DATE) {
                $tr.Text = $NEW_DATE
            }
        }
    }
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster

# Slide Master's own Date Placeholder shape.
Update-DatePlaceholder $master.Shapes

# Every Slide Layout's Date Placeholder shape.
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

Write-Host "Updated datetimeFigureOut placeholders on master + $($layouts.Count) layouts"
